# Applies "added ifo gdp component analysis preprocessing" edit:
# Extends the diagonal staircase of AVERAGE values in rows 45-53 by one
# more column to the right, updating the previously-last value in each
# row and filling in the newly-added trailing value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: update I45, add J45
$ws.Range("I45").Value = 0.4487415504340581
$ws.Range("J45").Value = 0.2388379152847414

# Row 46: update H46, add I46
$ws.Range("H46").Value = 0.5843816406042994
$ws.Range("I46").Value = 0.3744780054549828

# Row 47: update G47, add H47
$ws.Range("G47").Value = 0.3435754587486348
$ws.Range("H47").Value = 0.1336718235993181

# Row 48: update F48, add G48
$ws.Range("F48").Value = 0.2982442434965384
$ws.Range("G48").Value = 0.08834060834722172

# Row 49: update E49, add F49
$ws.Range("E49").Value = 0.2313828215604846
$ws.Range("F49").Value = 0.02147918641116785

# Row 50: update D50, add E50
$ws.Range("D50").Value = 0.201796619203768
$ws.Range("E50").Value = -0.00810701594554874

# Row 51: update C51, add D51
$ws.Range("C51").Value = 0.1836459624741271
$ws.Range("D51").Value = -0.02625767267518964

# Row 52: update B52, add C52
$ws.Range("B52").Value = 0.1656141382254278
$ws.Range("C52").Value = -0.04428949692388896

# Row 53: add B53
$ws.Range("B53").Value = -0.09587373626955231
